$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "28.152.87"
$ws.Range("E2").Value = "  -0.65%  "

# Row 3
$ws.Range("D3").Value = "1.906.34"
$ws.Range("E3").Value = "  +0.63%  "

# Row 4
$ws.Range("E4").Value = "  +0.30%  "

# Row 5
$ws.Range("D5").Value = "'313.81"
$ws.Range("E5").Value = "  +0.95%  "

# Row 6
$ws.Range("E6").Value = "  +0.42%  "

# Row 7
$ws.Range("D7").Value = "'0.5029"
$ws.Range("E7").Value = "  -0.04%  "

# Row 8
$ws.Range("D8").Value = "'0.3906"
$ws.Range("E8").Value = "  -0.19%  "

# Row 9
$ws.Range("D9").Value = "'0.09187"
$ws.Range("E9").Value = "  -2.96%  "

# Row 10
$ws.Range("D10").Value = "'1.137"
$ws.Range("E10").Value = "  -0.83%  "

# Row 11
$ws.Range("D11").Value = "'41.99"
$ws.Range("E11").Value = "  +2.06%  "

# Row 12
$ws.Range("D12").Value = "'6.397"
$ws.Range("E12").Value = "  -1.40%  "

# Row 13
$ws.Range("D13").Value = "'20.90"
$ws.Range("E13").Value = "  -0.68%  "

# Row 14
$ws.Range("D14").Value = "1.913.88"
$ws.Range("E14").Value = "  +1.47%  "

# Row 15
$ws.Range("D15").Value = "'7.309"
$ws.Range("E15").Value = "  -1.61%  "

# Row 16
$ws.Range("D16").Value = "'1.002"
$ws.Range("E16").Value = "  +0.40%  "

# Row 17
$ws.Range("D17").Value = "'0.00001116"
$ws.Range("E17").Value = "  -0.99%  "

# Row 18
$ws.Range("D18").Value = "'92.52"
$ws.Range("E18").Value = "  -0.15%  "

# Row 19
$ws.Range("D19").Value = "'0.06640"
$ws.Range("E19").Value = "  +0.77%  "

# Row 20
$ws.Range("D20").Value = "'17.96"
$ws.Range("E20").Value = "  +0.78%  "

# Row 21
$ws.Range("D21").Value = "'1.001"
$ws.Range("E21").Value = "  +0.18%  "

# Row 22
$ws.Range("D22").Value = "'6.223"
$ws.Range("E22").Value = "  +0.09%  "

# Row 23
$ws.Range("D23").Value = "28.212.51"
$ws.Range("E23").Value = "  -0.68%  "

# Row 24
$ws.Range("D24").Value = "'11.48"
$ws.Range("E24").Value = "  +1.62%  "

# Row 25
$ws.Range("D25").Value = "'2.317"
$ws.Range("E25").Value = "  +1.92%  "

# Row 26
$ws.Range("B26").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C26").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D26").Value = "2.125.05"
$ws.Range("E26").Value = "  +1.13%  "

# Row 27
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").Value = "'2.568"
$ws.Range("E27").Value = "  -1.99%  "

# Row 28
$ws.Range("D28").Value = "'21.07"
$ws.Range("E28").Value = "  -1.09%  "

# Row 29
$ws.Range("D29").Value = "'158.12"
$ws.Range("E29").Value = "  -0.21%  "

# Row 30
$ws.Range("D30").Value = "'126.91"
$ws.Range("E30").Value = "  -0.79%  "

# Row 31
$ws.Range("D31").Value = "'1.091"
$ws.Range("E31").Value = "  +1.22%  "

# Row 32
$ws.Range("D32").Value = "'0.1069"
$ws.Range("E32").Value = "  +0.95%  "

# Row 33
$ws.Range("D33").Value = "'5.647"
$ws.Range("E33").Value = "  +0.07%  "

# Row 34
$ws.Range("D34").Value = "'3.623"
$ws.Range("E34").Value = "  +0.27%  "

# Row 35
$ws.Range("D35").Value = "'9.762"
$ws.Range("E35").Value = "  +3.20%  "

# Row 36
$ws.Range("D36").Value = "'0.06658"
$ws.Range("E36").Value = "  -1.38%  "

# Row 37
$ws.Range("D37").Value = "'0.02426"
$ws.Range("E37").Value = "  +0.18%  "

# Row 38
$ws.Range("D38").Value = "'0.2217"
$ws.Range("E38").Value = "  +1.25%  "

# Row 39
$ws.Range("D39").Value = "'1.232"
$ws.Range("E39").Value = "  -1.50%  "

# Row 40
$ws.Range("B40").Value = "TheSandbox"
$ws.Range("C40").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").Value = "'0.6653"
$ws.Range("E40").Value = "  +4.29%  "

# Row 41
$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "'1.269"
$ws.Range("E41").Value = "  +6.77%  "

# Row 42
$ws.Range("D42").Value = "'11.47"
$ws.Range("E42").Value = "  -0.40%  "

# Row 43
$ws.Range("D43").Value = "'4.994"
$ws.Range("E43").Value = "  -0.87%  "

# Row 44
$ws.Range("D44").Value = "'1.001"
$ws.Range("E44").Value = "  +0.34%  "

# Row 45
$ws.Range("D45").Value = "'0.6196"
$ws.Range("E45").Value = "  +3.17%  "

# Row 46
$ws.Range("E46").Value = "  -1.67%  "

# Row 47
$ws.Range("B47").Value = "WEMIXTOKEN"
$ws.Range("C47").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D47").Value = "'1.297"
$ws.Range("E47").Value = "  +1.63%  "

# Row 48
$ws.Range("B48").Value = "PancakeSwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D48").Value = "'3.701"
$ws.Range("E48").Value = "  +1.35%  "

# Row 49
$ws.Range("D49").Value = "'2.019"
$ws.Range("E49").Value = "  +0.57%  "

# Row 50
$ws.Range("D50").Value = "'121.78"
$ws.Range("E50").Value = "  -1.71%  "

# Row 51
$ws.Range("D51").Value = "'1.189"
$ws.Range("E51").Value = "  -0.91%  "
